$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.050144795270566
$ws.Cells.Item(2, 4).Value = 1.057177983285805
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.064702124603033
$ws.Cells.Item(2, 9).Value = 1.048915766747342
$ws.Cells.Item(2, 10).Value = 1.055179501004725
$ws.Cells.Item(2, 11).Value = 1.059913665615101
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.067417349171799
$ws.Cells.Item(2, 14).Value = 1.022154257867351

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.051201983916221
$ws.Cells.Item(3, 4).Value = 1.058053607784603
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.065758795710461
$ws.Cells.Item(3, 9).Value = 1.049263441346406
$ws.Cells.Item(3, 10).Value = 1.055885727573932
$ws.Cells.Item(3, 11).Value = 1.060602917881396
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.068288672234229
$ws.Cells.Item(3, 14).Value = 1.022394941985281

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.051885866789486
$ws.Cells.Item(4, 4).Value = 1.058620003244672
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.066442759436675
$ws.Cells.Item(4, 9).Value = 1.049486983013634
$ws.Cells.Item(4, 10).Value = 1.056341921089158
$ws.Cells.Item(4, 11).Value = 1.06104807873897
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.068852099296417
$ws.Cells.Item(4, 14).Value = 1.022550282735972

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.05217332636557
$ws.Cells.Item(5, 4).Value = 1.058858069695607
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.066730352073675
$ws.Cells.Item(5, 9).Value = 1.049580618286837
$ws.Cells.Item(5, 10).Value = 1.056533517606057
$ws.Cells.Item(5, 11).Value = 1.061235025275306
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.069088873709596
$ws.Cells.Item(5, 14).Value = 1.022615492664159

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.052221589479754
$ws.Cells.Item(6, 4).Value = 1.058898039385171
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.066778643325573
$ws.Cells.Item(6, 9).Value = 1.049596320026911
$ws.Cells.Item(6, 10).Value = 1.05656567652602
$ws.Cells.Item(6, 11).Value = 1.061266402745626
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.069128623878341
$ws.Cells.Item(6, 14).Value = 1.02262643610681

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.051889708013831
$ws.Cells.Item(7, 4).Value = 1.058623184482031
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.066446602050287
$ws.Cells.Item(7, 9).Value = 1.04948823551517
$ws.Cells.Item(7, 10).Value = 1.056344481946847
$ws.Cells.Item(7, 11).Value = 1.061050577510068
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.068855263442644
$ws.Cells.Item(7, 14).Value = 1.022551154448839

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.050502116490944
$ws.Cells.Item(8, 4).Value = 1.057473944349079
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.065059184563625
$ws.Cells.Item(8, 9).Value = 1.049033560482141
$ws.Cells.Item(8, 10).Value = 1.055418335589444
$ws.Cells.Item(8, 11).Value = 1.060146773780228
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.067711894977955
$ws.Cells.Item(8, 14).Value = 1.022235680564294

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.048055536547168
$ws.Cells.Item(9, 4).Value = 1.055447373112635
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.062616115616497
$ws.Cells.Item(9, 9).Value = 1.048221435426098
$ws.Cells.Item(9, 10).Value = 1.053780353271788
$ws.Cells.Item(9, 11).Value = 1.058547789368723
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.065694244674441
$ws.Cells.Item(9, 14).Value = 1.021676729322043

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.046423466579713
$ws.Cells.Item(10, 4).Value = 1.054095350129147
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.060988568694166
$ws.Cells.Item(10, 9).Value = 1.047672668468442
$ws.Cells.Item(10, 10).Value = 1.052684328683333
$ws.Cells.Item(10, 11).Value = 1.057477518005638
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.064347199167559
$ws.Cells.Item(10, 14).Value = 1.021302049057362

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.045716513338072
$ws.Cells.Item(11, 4).Value = 1.053509678714815
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.060284097721716
$ws.Cells.Item(11, 9).Value = 1.047433301970376
$ws.Cells.Item(11, 10).Value = 1.05220877834077
$ws.Cells.Item(11, 11).Value = 1.057013062019431
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.063763450585401
$ws.Cells.Item(11, 14).Value = 1.021139323149687

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.04545388015437
$ws.Cells.Item(12, 4).Value = 1.053292098693981
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.060022465774871
$ws.Cells.Item(12, 9).Value = 1.047344127998663
$ws.Cells.Item(12, 10).Value = 1.052031992576221
$ws.Cells.Item(12, 11).Value = 1.056840388729835
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.063546549659203
$ws.Cells.Item(12, 14).Value = 1.021078806321854

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.045510217653585
$ws.Cells.Item(13, 4).Value = 1.053338771957857
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.06007858491252
$ws.Cells.Item(13, 9).Value = 1.047363268003952
$ws.Cells.Item(13, 10).Value = 1.052069920297869
$ws.Cells.Item(13, 11).Value = 1.056877434693097
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.063593078848049
$ws.Cells.Item(13, 14).Value = 1.021091790698935

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.045694804802869
$ws.Cells.Item(14, 4).Value = 1.053491694208509
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.060262470330158
$ws.Cells.Item(14, 9).Value = 1.047425936179882
$ws.Cells.Item(14, 10).Value = 1.052194168137929
$ws.Cells.Item(14, 11).Value = 1.056998791927149
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.063745522931128
$ws.Cells.Item(14, 14).Value = 1.021134322302863

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.045808529810199
$ws.Cells.Item(15, 4).Value = 1.053585909938113
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.060375773489359
$ws.Cells.Item(15, 9).Value = 1.047464513309953
$ws.Cells.Item(15, 10).Value = 1.052270702084546
$ws.Cells.Item(15, 11).Value = 1.0570735437574
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.063839439384526
$ws.Cells.Item(15, 14).Value = 1.021160517729699

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.046470379273719
$ws.Cells.Item(16, 4).Value = 1.054134214224808
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.061035327727808
$ws.Cells.Item(16, 9).Value = 1.047688517618023
$ws.Cells.Item(16, 10).Value = 1.0527158690482
$ws.Cells.Item(16, 11).Value = 1.057508320882614
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.064385930735972
$ws.Cells.Item(16, 14).Value = 1.021312838379181

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.046885470816928
$ws.Cells.Item(17, 4).Value = 1.05447808765996
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.061449120131841
$ws.Cells.Item(17, 9).Value = 1.047828561824449
$ws.Cells.Item(17, 10).Value = 1.052994852102362
$ws.Cells.Item(17, 11).Value = 1.057780771482699
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.064728604584494
$ws.Cells.Item(17, 14).Value = 1.02140825472405

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.047127562026757
$ws.Cells.Item(18, 4).Value = 1.054678640397146
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.061690504076137
$ws.Cells.Item(18, 9).Value = 1.047910078712234
$ws.Cells.Item(18, 10).Value = 1.053157485160917
$ws.Cells.Item(18, 11).Value = 1.057939588769153
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.064928435209118
$ws.Cells.Item(18, 14).Value = 1.021463862488947

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.04721010471345
$ws.Cells.Item(19, 4).Value = 1.054747019828428
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.061772814153385
$ws.Cells.Item(19, 9).Value = 1.047937845283958
$ws.Cells.Item(19, 10).Value = 1.053212923059308
$ws.Cells.Item(19, 11).Value = 1.05799372465172
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.064996564623229
$ws.Cells.Item(19, 14).Value = 1.021482815332982

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.046840937979394
$ws.Cells.Item(20, 4).Value = 1.054441195655284
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.06140472145241
$ws.Cells.Item(20, 9).Value = 1.047813553839247
$ws.Cells.Item(20, 10).Value = 1.052964929486165
$ws.Cells.Item(20, 11).Value = 1.057751550305437
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.064691843603987
$ws.Cells.Item(20, 14).Value = 1.02139802231227

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.045640449578542
$ws.Cells.Item(21, 4).Value = 1.053446663399733
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.060208319557496
$ws.Cells.Item(21, 9).Value = 1.047407489214283
$ws.Cells.Item(21, 10).Value = 1.05215758425624
$ws.Cells.Item(21, 11).Value = 1.056963059489745
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.063700633902233
$ws.Cells.Item(21, 14).Value = 1.021121799825972

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.04488542726557
$ws.Cells.Item(22, 4).Value = 1.052821155616698
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.059456325144076
$ws.Cells.Item(22, 9).Value = 1.047150660885262
$ws.Cells.Item(22, 10).Value = 1.051649134076204
$ws.Cells.Item(22, 11).Value = 1.056466414910661
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.063077011782234
$ws.Cells.Item(22, 14).Value = 1.020947704370785

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.045285700716057
$ws.Cells.Item(23, 4).Value = 1.053152768657738
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.059854949808335
$ws.Cells.Item(23, 9).Value = 1.047286954515853
$ws.Cells.Item(23, 10).Value = 1.051918752875603
$ws.Cells.Item(23, 11).Value = 1.05672977994295
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.063407644542844
$ws.Cells.Item(23, 14).Value = 1.021040035811134

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.0468610605277
$ws.Cells.Item(24, 4).Value = 1.054457865634568
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.061424783224505
$ws.Cells.Item(24, 9).Value = 1.047820335822362
$ws.Cells.Item(24, 10).Value = 1.052978450515716
$ws.Cells.Item(24, 11).Value = 1.057764754401943
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.064708454448304
$ws.Cells.Item(24, 14).Value = 1.021402646043802

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.048688213046016
$ws.Cells.Item(25, 4).Value = 1.05597146288888
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.063247501822515
$ws.Cells.Item(25, 9).Value = 1.048432684615081
$ws.Cells.Item(25, 10).Value = 1.054204521935244
$ws.Cells.Item(25, 11).Value = 1.058961919765874
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.066216198291762
$ws.Cells.Item(25, 14).Value = 1.021821592388107

